$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.7047556261015829
$ws.Range("D2").Value = -0.1005524227762842
$ws.Range("E2").Value = -0.04743437747629582
$ws.Range("F2").Value = 0.02275530254629318
$ws.Range("H2").Value = 0.01984950350403698
$ws.Range("I2").Value = -0.07735816559944704
$ws.Range("J2").Value = -0.1403061574535686
$ws.Range("K2").Value = -0.04671663250852448
$ws.Range("L2").Value = 0.02676678503401862
$ws.Range("M2").Value = -0.07828834604569371
$ws.Range("N2").Value = 0.004732017893733165
$ws.Range("O2").Value = -0.08118141911519945
$ws.Range("P2").Value = -0.02673676715644199
$ws.Range("B3").Value = -0.7047556261015829
$ws.Range("D3").Value = 0.140323046109499
$ws.Range("B4").Value = -0.1005524227762842
$ws.Range("C4").Value = 0.140323046109499
$ws.Range("E4").Value = 0.2912621032440872
$ws.Range("F4").Value = -0.1478776120475858
$ws.Range("H4").Value = -0.1909117508911267
$ws.Range("I4").Value = 0.3510319043220564
$ws.Range("J4").Value = -0.4814546273730555
$ws.Range("K4").Value = 0.03448084153235571
$ws.Range("L4").Value = -0.3926606960539389
$ws.Range("M4").Value = 0.0682481025881473
$ws.Range("N4").Value = -0.2865747988606426
$ws.Range("O4").Value = 0.2750477972531822
$ws.Range("P4").Value = -0.4368722050314565
$ws.Range("B5").Value = -0.04743437747629582
$ws.Range("D5").Value = 0.2912621032440872
$ws.Range("B6").Value = 0.02275530254629318
$ws.Range("D6").Value = -0.1478776120475858
$ws.Range("B8").Value = 0.01984950350403698
$ws.Range("D8").Value = -0.1909117508911267
$ws.Range("B9").Value = -0.07735816559944704
$ws.Range("D9").Value = 0.3510319043220564
$ws.Range("B10").Value = -0.1403061574535686
$ws.Range("D10").Value = -0.4814546273730555
$ws.Range("B11").Value = -0.04671663250852448
$ws.Range("D11").Value = 0.03448084153235571
$ws.Range("B12").Value = 0.02676678503401862
$ws.Range("D12").Value = -0.3926606960539389
$ws.Range("B13").Value = -0.07828834604569371
$ws.Range("D13").Value = 0.0682481025881473
$ws.Range("B14").Value = 0.004732017893733165
$ws.Range("D14").Value = -0.2865747988606426
$ws.Range("B15").Value = -0.08118141911519945
$ws.Range("D15").Value = 0.2750477972531822
$ws.Range("B16").Value = -0.02673676715644199
$ws.Range("D16").Value = -0.4368722050314565
